# PAMGuard_exporter_dialog_annotated.pptx -- "X3 fix and building for MACOS"
#
# 1) Text correction on the slide: "e.g." -> "e.g.," in the
#    "Data specific settings ..." callout.
# 2) Refresh of the "update automatically" Date placeholder fields
#    (slide master + every slide layout) that PowerPoint re-stamped
#    the next time the deck was saved (8/22/24 -> 10/7/24).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Fix the wording on the slide itself.
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Data specific settings for filtering e.g. selecting on classified clicks") {
            # Remember the shape's current size: it auto-fits its text box,
            # and re-assigning .Text recalculates the fit. Restore the
            # original height afterwards so only the wording changes.
            $origHeight = $shp.Height
            $origWidth  = $shp.Width
            $tr.Text = "Data specific settings for filtering e.g., selecting on classified clicks"
            $shp.Height = $origHeight
            $shp.Width  = $origWidth
        }
    }
}

# ---------------------------------------------------------------------
# 2) Refresh the "Date automatically" placeholders site-wide.
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes, $oldText, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $tr = $shp.TextFrame.TextRange
                if ($tr.Text -eq $oldText) {
                    $tr.Text = $newText
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes "8/22/24" "10/7/24"

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes "8/22/24" "10/7/24"
}
